{"js": "// Replace the date and each division-problem cell with the new values,\n// matching the original text to the corresponding new text in document\n// (reading) order. Using positional/ordered matching avoids accidentally\n// re-touching text that a previous replacement happens to have produced.\nconst replacements = [\n  [\"2024-05-16 Thursday\", \"2024-05-17 Friday\"],\n  [\"86\u00f79=\", \"63\u00f74=\"],\n  [\"87\u00f72=\", \"70\u00f74=\"],\n  [\"45\u00f75=\", \"95\u00f73=\"],\n  [\"49\u00f77=\", \"19\u00f72=\"],\n  [\"97\u00f78=\", \"65\u00f73=\"],\n  [\"33\u00f75=\", \"28\u00f76=\"],\n  [\"87\u00f77=\", \"99\u00f79=\"],\n  [\"92\u00f77=\", \"97\u00f76=\"],\n  [\"52\u00f79=\", \"74\u00f76=\"],\n  [\"64\u00f76=\", \"74\u00f77=\"],\n  [\"48\u00f75=\", \"21\u00f78=\"],\n  [\"99\u00f79=\", \"79\u00f72=\"],\n  [\"18\u00f74=\", \"57\u00f78=\"],\n  [\"30\u00f74=\", \"13\u00f79=\"],\n  [\"26\u00f73=\", \"37\u00f77=\"],\n  [\"40\u00f79=\", \"57\u00f75=\"],\n  [\"25\u00f74=\", \"15\u00f77=\"],\n  [\"10\u00f72=\", \"92\u00f77=\"],\n  [\"62\u00f72=\", \"72\u00f78=\"],\n  [\"82\u00f75=\", \"79\u00f77=\"],\n  [\"59\u00f77=\", \"35\u00f72=\"],\n  [\"98\u00f73=\", \"88\u00f72=\"],\n  [\"80\u00f74=\", \"71\u00f72=\"],\n  [\"69\u00f74=\", \"95\u00f74=\"],\n  [\"75\u00f74=\", \"44\u00f78=\"],\n];\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet idx = 0;\nfor (let i = 0; i < paragraphs.items.length && idx < replacements.length; i++) {\n  const para = paragraphs.items[i];\n  const [oldText, newText] = replacements[idx];\n  if (para.text === oldText) {\n    para.insertText(newText, Word.InsertLocation.replace);\n    idx++;\n  }\n}\n\nawait context.sync();\n\nif (idx !== replacements.length) {\n  throw new Error(\n    \"Only matched \" + idx + \" of \" + replacements.length + \" expected paragraphs\"\n  );\n}\n", "ps1": "# Replace the date and each division-problem cell with the new values,\n# matching the original text to the corresponding new text in document\n# (reading) order. Positional/ordered matching avoids accidentally\n# re-touching text that a previous replacement happens to have produced\n# (a few new values coincide with other cells' old values).\n$replacements = @(\n    @{Old=\"2024-05-16 Thursday\"; New=\"2024-05-17 Friday\"};\n    @{Old=\"86\u00f79=\"; New=\"63\u00f74=\"};\n    @{Old=\"87\u00f72=\"; New=\"70\u00f74=\"};\n    @{Old=\"45\u00f75=\"; New=\"95\u00f73=\"};\n    @{Old=\"49\u00f77=\"; New=\"19\u00f72=\"};\n    @{Old=\"97\u00f78=\"; New=\"65\u00f73=\"};\n    @{Old=\"33\u00f75=\"; New=\"28\u00f76=\"};\n    @{Old=\"87\u00f77=\"; New=\"99\u00f79=\"};\n    @{Old=\"92\u00f77=\"; New=\"97\u00f76=\"};\n    @{Old=\"52\u00f79=\"; New=\"74\u00f76=\"};\n    @{Old=\"64\u00f76=\"; New=\"74\u00f77=\"};\n    @{Old=\"48\u00f75=\"; New=\"21\u00f78=\"};\n    @{Old=\"99\u00f79=\"; New=\"79\u00f72=\"};\n    @{Old=\"18\u00f74=\"; New=\"57\u00f78=\"};\n    @{Old=\"30\u00f74=\"; New=\"13\u00f79=\"};\n    @{Old=\"26\u00f73=\"; New=\"37\u00f77=\"};\n    @{Old=\"40\u00f79=\"; New=\"57\u00f75=\"};\n    @{Old=\"25\u00f74=\"; New=\"15\u00f77=\"};\n    @{Old=\"10\u00f72=\"; New=\"92\u00f77=\"};\n    @{Old=\"62\u00f72=\"; New=\"72\u00f78=\"};\n    @{Old=\"82\u00f75=\"; New=\"79\u00f77=\"};\n    @{Old=\"59\u00f77=\"; New=\"35\u00f72=\"};\n    @{Old=\"98\u00f73=\"; New=\"88\u00f72=\"};\n    @{Old=\"80\u00f74=\"; New=\"71\u00f72=\"};\n    @{Old=\"69\u00f74=\"; New=\"95\u00f74=\"};\n    @{Old=\"75\u00f74=\"; New=\"44\u00f78=\"};\n)\n\n$d = $word.ActiveDocument\n\n$idx = 0\nforeach ($p in $d.Paragraphs) {\n    if ($idx -ge $replacements.Count) { break }\n    $rng = $p.Range\n    $core = $rng.Text.TrimEnd([char]13, [char]7)\n    $pair = $replacements[$idx]\n    if ($core -eq $pair.Old) {\n        $rng.Text = $pair.New\n        $idx = $idx + 1\n    }\n}\n\nif ($idx -ne $replacements.Count) {\n    throw \"Only matched $idx of $($replacements.Count) expected paragraphs\"\n}\n"}
